{"js": "// Replace the division expressions found in the table cells with their\n// updated values, in document order. Each \"old\" text is unique in the\n// document, so a single targeted search-and-replace per pair is safe even\n// though some replacement values momentarily duplicate other original\n// values (e.g. \"17\u00f78=\" appears both as an original value earlier in the\n// document and as a replacement value later on).\nconst replacements = [\n  [\"77\u00f77=\", \"73\u00f72=\"],\n  [\"34\u00f79=\", \"58\u00f74=\"],\n  [\"38\u00f78=\", \"45\u00f75=\"],\n  [\"74\u00f74=\", \"35\u00f77=\"],\n  [\"74\u00f79=\", \"51\u00f78=\"],\n  [\"72\u00f79=\", \"70\u00f76=\"],\n  [\"61\u00f76=\", \"47\u00f74=\"],\n  [\"56\u00f74=\", \"29\u00f79=\"],\n  [\"71\u00f75=\", \"97\u00f76=\"],\n  [\"35\u00f74=\", \"51\u00f78=\"],\n  [\"52\u00f73=\", \"84\u00f78=\"],\n  [\"60\u00f75=\", \"74\u00f76=\"],\n  [\"64\u00f75=\", \"90\u00f74=\"],\n  [\"44\u00f75=\", \"46\u00f73=\"],\n  [\"89\u00f72=\", \"21\u00f77=\"],\n  [\"14\u00f74=\", \"74\u00f77=\"],\n  [\"97\u00f73=\", \"92\u00f73=\"],\n  [\"17\u00f78=\", \"13\u00f76=\"],\n  [\"33\u00f77=\", \"31\u00f75=\"],\n  [\"88\u00f78=\", \"61\u00f77=\"],\n  [\"59\u00f78=\", \"17\u00f78=\"],\n  [\"48\u00f76=\", \"11\u00f79=\"],\n  [\"96\u00f74=\", \"64\u00f78=\"],\n  [\"20\u00f77=\", \"95\u00f77=\"],\n  [\"65\u00f79=\", \"36\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Only the first match is relevant: by construction every \"old\" value is\n  // unique in the document at the time it is searched for (replacement\n  // values that coincide with a not-yet-processed original value only\n  // appear later in document order, after the original has already been\n  // handled).\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the division expressions found in the table cells with their\n# updated values, in document order. Each \"old\" text is unique in the\n# document, so a single Find/Replace (first match, not ReplaceAll) per pair\n# is safe even though some replacement values momentarily duplicate other\n# original values later in the document (e.g. \"17\u00f78=\" appears both as an\n# original value earlier in the document and as a replacement value later\n# on) -- processing strictly top-to-bottom avoids any cross-talk.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"77\u00f77=\", \"73\u00f72=\"),\n    @(\"34\u00f79=\", \"58\u00f74=\"),\n    @(\"38\u00f78=\", \"45\u00f75=\"),\n    @(\"74\u00f74=\", \"35\u00f77=\"),\n    @(\"74\u00f79=\", \"51\u00f78=\"),\n    @(\"72\u00f79=\", \"70\u00f76=\"),\n    @(\"61\u00f76=\", \"47\u00f74=\"),\n    @(\"56\u00f74=\", \"29\u00f79=\"),\n    @(\"71\u00f75=\", \"97\u00f76=\"),\n    @(\"35\u00f74=\", \"51\u00f78=\"),\n    @(\"52\u00f73=\", \"84\u00f78=\"),\n    @(\"60\u00f75=\", \"74\u00f76=\"),\n    @(\"64\u00f75=\", \"90\u00f74=\"),\n    @(\"44\u00f75=\", \"46\u00f73=\"),\n    @(\"89\u00f72=\", \"21\u00f77=\"),\n    @(\"14\u00f74=\", \"74\u00f77=\"),\n    @(\"97\u00f73=\", \"92\u00f73=\"),\n    @(\"17\u00f78=\", \"13\u00f76=\"),\n    @(\"33\u00f77=\", \"31\u00f75=\"),\n    @(\"88\u00f78=\", \"61\u00f77=\"),\n    @(\"59\u00f78=\", \"17\u00f78=\"),\n    @(\"48\u00f76=\", \"11\u00f79=\"),\n    @(\"96\u00f74=\", \"64\u00f78=\"),\n    @(\"20\u00f77=\", \"95\u00f77=\"),\n    @(\"65\u00f79=\", \"36\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
